$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-15 Sunday" "2026-02-16 Monday"

Replace-Text "33÷5=" "95÷8="
Replace-Text "58÷5=" "53÷7="
Replace-Text "80÷8=" "65÷9="
Replace-Text "98÷8=" "84÷9="
Replace-Text "34÷2=" "28÷2="

Replace-Text "47÷2=" "55÷2="
Replace-Text "48÷6=" "20÷8="
Replace-Text "96÷9=" "60÷2="
Replace-Text "58÷9=" "92÷3="

Replace-Text "74÷3=" "58÷4="
Replace-Text "73÷4=" "38÷4="
Replace-Text "48÷5=" "15÷9="
Replace-Text "56÷7=" "22÷4="
Replace-Text "66÷4=" "62÷2="

Replace-Text "14÷8=" "23÷6="
Replace-Text "44÷7=" "95÷3="
Replace-Text "96÷2=" "68÷7="
Replace-Text "76÷2=" "84÷3="
Replace-Text "67÷4=" "76÷6="

Replace-Text "21÷2=" "36÷4="
Replace-Text "82÷6=" "56÷8="
Replace-Text "40÷8=" "91÷3="
Replace-Text "89÷6=" "98÷5="
Replace-Text "51÷9=" "38÷3="
